# Update of Excel Modules Files
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "ElementName3" header (column F) and its data cell.
$ws.Range("F1").ClearContents()

# Fill in the new module / element data rows.
$ws.Range("A2").Value = 'GSTR51'
$ws.Range("B2").Value = 'pede. Suspendisse dui.'
$ws.Range("C2").Value = 'EL Haddad'
$ws.Range("D2").Value = 'Nullam feugiat placerat'
$ws.Range("E2").Value = 'varius et, euismod'

$ws.Range("A3").Value = 'GSTR52'
$ws.Range("B3").Value = 'a nunc. In'
$ws.Range("C3").Value = 'Badir'
$ws.Range("D3").Value = 'sodales nisi magna'
$ws.Range("E3").Value = 'elementum sem, vitae'

$ws.Range("A4").Value = 'GSTR53'
$ws.Range("B4").Value = 'amet metus. Aliquam'
$ws.Range("C4").Value = 'Ezzine'
$ws.Range("D4").Value = 'Cras vulputate velit'
$ws.Range("E4").Value = 'scelerisque neque sed'

$ws.Range("A5").Value = 'GSTR54'
$ws.Range("B5").Value = 'quam vel sapien'
$ws.Range("C5").Value = 'El Alami Hassoun'
$ws.Range("D5").Value = 'Nunc mauris elit,'
$ws.Range("E5").Value = 'libero et tristique'

$ws.Range("A6").Value = 'GSTR55'
$ws.Range("B6").Value = 'feugiat nec, diam.'
$ws.Range("C6").Value = 'Lazaar'
$ws.Range("D6").Value = 'pellentesque. Sed dictum.'
$ws.Range("E6").Value = 'ridiculus mus. Proin'

$ws.Range("A7").Value = 'GSTR54-2'
$ws.Range("B7").Value = 'nonummy. Fusce fermentum'
$ws.Range("C7").Value = 'El Haddad'
$ws.Range("D7").Value = 'neque pellentesque massa'
$ws.Range("E7").Value = 'Mauris eu turpis.'

$ws.Range("A8").Value = 'GSTR55-2'
$ws.Range("B8").Value = 'a, arcu. Sed'
$ws.Range("C8").Value = 'EL Haddad'
$ws.Range("D8").Value = 'sit amet risus.'
$ws.Range("E8").Value = 'Nulla facilisi. Sed'

$ws.Range("A9").Value = 'GSTR56'
$ws.Range("B9").Value = 'Suspendisse eleifend. Cras'
$ws.Range("C9").Value = 'El Alami Hassoun'
$ws.Range("D9").Value = 'velit dui, semper'
$ws.Range("E9").Value = 'ligula elit, pretium'

# B2 carries explicit black font color formatting.
$ws.Range("B2").Font.Color = 0

# Resize the data columns to fit the new, wider content.
$ws.Columns.Item(1).ColumnWidth = 8.333333333333334
$ws.Columns.Item(2).ColumnWidth = 26.166666666666668
$ws.Columns.Item(3).ColumnWidth = 15.333333333333334
$ws.Columns.Item(4).ColumnWidth = 24.333333333333332
$ws.Columns.Item(5).ColumnWidth = 20.5

# Restore the previous selection location.
$ws.Range("H12").Select() | Out-Null
